# ---------------------------------------------------------------------------
# daily_submission_stats.xlsx update
#   - turn the numeric task-id column (A2:A4) into the task-name labels
#   - append three more "score > N" mini-tables (rows 7-26) with their own
#     header/footer rows copied (styles + shared strings) from the existing
#     table at the top of the sheet
#   - rewire the existing chart series onto the new label column + add
#     3 new series for the "score > 0" block
#   - reposition the chart to the right of the new tables
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. relabel the first table's task column -------------------------------
$ws.Range("A2").Value = "rect"
$ws.Range("A3").Value = "shoes"
$ws.Range("A4").Value = "split"

# --- 2. "score > 0" block (rows 7-12) ---------------------------------------
$ws.Range("A7").Value = "score > 0"

$ws.Range("A1:K1").Copy($ws.Range("A8:K8"))

$ws.Range("A2:K2").Copy($ws.Range("A9:K9"))
$ws.Range("A9").Value = "rect"
$ws.Range("B9").Value = 11
$ws.Range("C9").Value = 12
$ws.Range("D9").Value = 47
$ws.Range("E9").Value = 72
$ws.Range("F9").Value = 87
$ws.Range("G9").Value = 101
$ws.Range("H9").Value = 108
$ws.Range("I9").Value = 143
$ws.Range("J9").Value = 179
$ws.Range("K9").Value = 308

$ws.Range("A2:K2").Copy($ws.Range("A10:K10"))
$ws.Range("A10").Value = "shoes"
$ws.Range("B10").Value = 163
$ws.Range("C10").Value = 234
$ws.Range("D10").Value = 294
$ws.Range("E10").Value = 172
$ws.Range("F10").Value = 216
$ws.Range("G10").Value = 137
$ws.Range("H10").Value = 167
$ws.Range("I10").Value = 121
$ws.Range("J10").Value = 115
$ws.Range("K10").Value = 155

$ws.Range("A2:K2").Copy($ws.Range("A11:K11"))
$ws.Range("A11").Value = "split"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 18
$ws.Range("D11").Value = 62
$ws.Range("E11").Value = 95
$ws.Range("F11").Value = 144
$ws.Range("G11").Value = 148
$ws.Range("H11").Value = 165
$ws.Range("I11").Value = 211
$ws.Range("J11").Value = 164
$ws.Range("K11").Value = 287

$ws.Range("B5:K5").Copy($ws.Range("B12:K12"))

# --- 3. "score > 20" block (rows 14-19) -------------------------------------
$ws.Range("A14").Value = "score > 20"
$ws.Range("A14").WrapText = $false

$ws.Range("A1:K1").Copy($ws.Range("A15:K15"))

$ws.Range("A2:K2").Copy($ws.Range("A16:K16"))
$ws.Range("A16").Value = "rect"
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 19
$ws.Range("E16").Value = 52
$ws.Range("F16").Value = 58
$ws.Range("G16").Value = 62
$ws.Range("H16").Value = 73
$ws.Range("I16").Value = 84
$ws.Range("J16").Value = 116
$ws.Range("K16").Value = 190

$ws.Range("A2:K2").Copy($ws.Range("A17:K17"))
$ws.Range("A17").Value = "shoes"
$ws.Range("B17").Value = 55
$ws.Range("C17").Value = 107
$ws.Range("D17").Value = 169
$ws.Range("E17").Value = 100
$ws.Range("F17").Value = 114
$ws.Range("G17").Value = 63
$ws.Range("H17").Value = 102
$ws.Range("I17").Value = 62
$ws.Range("J17").Value = 57
$ws.Range("K17").Value = 77

$ws.Range("A2:K2").Copy($ws.Range("A18:K18"))
$ws.Range("A18").Value = "split"
$ws.Range("B18").Value = 4
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 15
$ws.Range("E18").Value = 23
$ws.Range("F18").Value = 64
$ws.Range("G18").Value = 65
$ws.Range("H18").Value = 66
$ws.Range("I18").Value = 99
$ws.Range("J18").Value = 62
$ws.Range("K18").Value = 122

$ws.Range("B5:K5").Copy($ws.Range("B19:K19"))

# --- 4. "score > 50" block (rows 21-26) -------------------------------------
$ws.Range("A21").Value = "score > 50"

$ws.Range("A1:K1").Copy($ws.Range("A22:K22"))

$ws.Range("A2:K2").Copy($ws.Range("A23:K23"))
$ws.Range("A23").Value = "rect"
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 0
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 24
$ws.Range("F23").Value = 21
$ws.Range("G23").Value = 23
$ws.Range("H23").Value = 28
$ws.Range("I23").Value = 30
$ws.Range("J23").Value = 31
$ws.Range("K23").Value = 58

$ws.Range("A2:K2").Copy($ws.Range("A24:K24"))
$ws.Range("A24").Value = "shoes"
$ws.Range("B24").Value = 4
$ws.Range("C24").Value = 40
$ws.Range("D24").Value = 52
$ws.Range("E24").Value = 27
$ws.Range("F24").Value = 27
$ws.Range("G24").Value = 15
$ws.Range("H24").Value = 26
$ws.Range("I24").Value = 11
$ws.Range("J24").Value = 13
$ws.Range("K24").Value = 8

$ws.Range("A2:K2").Copy($ws.Range("A25:K25"))
$ws.Range("A25").Value = "split"
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 2
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 8

$ws.Range("B5:K5").Copy($ws.Range("B26:K26"))

# --- 5. rewire chart series onto the new label column + new cat range ------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()

$s1 = $sc.Item(1)
$s1.Name = "rect"
$s1.Formula = "=SERIES(""rect"",'2019 DAY 1'!`$B`$12:`$K`$12,'2019 DAY 1'!`$B`$2:`$K`$2,1)"

$s2 = $sc.Item(2)
$s2.Name = "shoes"
$s2.Formula = "=SERIES(""shoes"",'2019 DAY 1'!`$B`$12:`$K`$12,'2019 DAY 1'!`$B`$3:`$K`$3,2)"

$s3 = $sc.Item(3)
$s3.Name = "split"
$s3.Formula = "=SERIES(""split"",'2019 DAY 1'!`$B`$12:`$K`$12,'2019 DAY 1'!`$B`$4:`$K`$4,3)"

# --- 6. add the 3 new "score > 0" series ------------------------------------
$s4 = $sc.NewSeries()
$s4.Name = "rect0"
$s4.Formula = "=SERIES(""rect0"",'2019 DAY 1'!`$B`$12:`$K`$12,'2019 DAY 1'!`$B`$9:`$K`$9,4)"

$s5 = $sc.NewSeries()
$s5.Name = "shoes0"
$s5.Formula = "=SERIES(""shoes0"",'2019 DAY 1'!`$B`$12:`$K`$12,'2019 DAY 1'!`$B`$10:`$K`$10,5)"

$s6 = $sc.NewSeries()
$s6.Name = "split0"
$s6.Formula = "=SERIES(""split0"",'2019 DAY 1'!`$B`$12:`$K`$12,'2019 DAY 1'!`$B`$11:`$K`$11,6)"

# --- 7. move the chart up and to the right of the new tables ----------------
$co.Left = 673.5625
$co.Top = 13.87496062992126
$co.Width = 443.5
$co.Height = 216.0

# --- 8. restore the selection shown in the saved file -----------------------
$ws.Range("Q18").Select()

Write-Host "edit complete"
